$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")
try {
  $ws.Rows(67).Insert()
  Write-Host "insert ok"
} catch {
  Write-Host "insert err: $_"
}
